$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and reporting week) ---
$ws.Range("A8").Value = "Volume 30   Number  29"
$ws.Range("C9").Value = "Report Covering the Week  7/17/2023  Through  7/23/2023"

# --- Weekly crime-stat numeric updates (rows 14-30) ---
$ws.Range("C14").Value = 2
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = -60
$ws.Range("F14").Value = 8
$ws.Range("G14").Value = 11
$ws.Range("H14").Value = -27.272727272727
$ws.Range("I14").Value = 34
$ws.Range("J14").Value = 29
$ws.Range("K14").Value = 17.241379310344
$ws.Range("L14").Value = -2.857142857142
$ws.Range("M14").Value = 17.241379310344
$ws.Range("N14").Value = -82.198952879581
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = 18
$ws.Range("H15").Value = -55.555555555555
$ws.Range("I15").Value = 74
$ws.Range("J15").Value = 110
$ws.Range("K15").Value = -32.727272727272
$ws.Range("L15").Value = -5.128205128205
$ws.Range("M15").Value = -32.727272727272
$ws.Range("N15").Value = -73.188405797101
$ws.Range("C16").Value = 37
$ws.Range("D16").Value = 51
$ws.Range("E16").Value = -27.450980392156
$ws.Range("F16").Value = 161
$ws.Range("G16").Value = 181
$ws.Range("H16").Value = -11.049723756906
$ws.Range("I16").Value = 1037
$ws.Range("J16").Value = 1142
$ws.Range("K16").Value = -9.194395796847
$ws.Range("L16").Value = 12.472885032538
$ws.Range("M16").Value = -18.474842767295
$ws.Range("N16").Value = -80.183451175234
$ws.Range("C17").Value = 69
$ws.Range("D17").Value = 71
$ws.Range("E17").Value = -2.81690140845
$ws.Range("F17").Value = 252
$ws.Range("G17").Value = 264
$ws.Range("H17").Value = -4.545454545454
$ws.Range("I17").Value = 1661
$ws.Range("J17").Value = 1727
$ws.Range("K17").Value = -3.821656050955
$ws.Range("L17").Value = 9.854497354497
$ws.Range("M17").Value = 53.796296296296
$ws.Range("N17").Value = -49.909529553679
$ws.Range("C18").Value = 26
$ws.Range("E18").Value = -27.777777777777
$ws.Range("F18").Value = 115
$ws.Range("H18").Value = -8
$ws.Range("I18").Value = 857
$ws.Range("J18").Value = 1021
$ws.Range("K18").Value = -16.062683643486
$ws.Range("L18").Value = 12.319790301441
$ws.Range("M18").Value = 16.124661246612
$ws.Range("N18").Value = -85.808908759728
$ws.Range("C19").Value = 126
$ws.Range("E19").Value = -8.695652173913
$ws.Range("F19").Value = 507
$ws.Range("G19").Value = 551
$ws.Range("H19").Value = -7.985480943738
$ws.Range("I19").Value = 3393
$ws.Range("J19").Value = 3570
$ws.Range("K19").Value = -4.957983193277
$ws.Range("L19").Value = 27.412692452121
$ws.Range("M19").Value = 33.320235756385
$ws.Range("N19").Value = -43.76864434869
$ws.Range("C20").Value = 31
$ws.Range("D20").Value = 27
$ws.Range("E20").Value = 14.814814814814
$ws.Range("F20").Value = 119
$ws.Range("G20").Value = 115
$ws.Range("H20").Value = 3.478260869565
$ws.Range("I20").Value = 766
$ws.Range("J20").Value = 722
$ws.Range("K20").Value = 6.094182825484
$ws.Range("L20").Value = 58.921161825726
$ws.Range("M20").Value = 149.511400651466
$ws.Range("N20").Value = -85.415079969535
$ws.Range("C21").Value = 293
$ws.Range("D21").Value = 335
$ws.Range("E21").Value = -12.537313432835
$ws.Range("F21").Value = 1170
$ws.Range("G21").Value = 1265
$ws.Range("H21").Value = -7.509881422924
$ws.Range("I21").Value = 7822
$ws.Range("J21").Value = 8321
$ws.Range("K21").Value = -5.996875375555
$ws.Range("L21").Value = 21.177381874515
$ws.Range("M21").Value = 28.630159513238
$ws.Range("N21").Value = -70.304847955658
$ws.Range("C22").Value = 3
$ws.Range("D22").Value = 5
$ws.Range("E22").Value = -40
$ws.Range("F22").Value = 18
$ws.Range("G22").Value = 25
$ws.Range("H22").Value = -28
$ws.Range("I22").Value = 159
$ws.Range("J22").Value = 177
$ws.Range("K22").Value = -10.169491525423
$ws.Range("L22").Value = 28.225806451612
$ws.Range("M22").Value = 21.374045801526
$ws.Range("C23").Value = 28
$ws.Range("D23").Value = 32
$ws.Range("E23").Value = -12.5
$ws.Range("F23").Value = 100
$ws.Range("G23").Value = 103
$ws.Range("H23").Value = -2.912621359223
$ws.Range("I23").Value = 686
$ws.Range("J23").Value = 718
$ws.Range("K23").Value = -4.456824512534
$ws.Range("L23").Value = 2.388059701492
$ws.Range("M23").Value = 50.76923076923
$ws.Range("C24").Value = 293
$ws.Range("D24").Value = 323
$ws.Range("E24").Value = -9.287925696594
$ws.Range("F24").Value = 1154
$ws.Range("G24").Value = 1271
$ws.Range("H24").Value = -9.205350118017
$ws.Range("I24").Value = 7885
$ws.Range("J24").Value = 8795
$ws.Range("K24").Value = -10.346787947697
$ws.Range("L24").Value = 19.216812821288
$ws.Range("M24").Value = 58.333333333333
$ws.Range("C25").Value = 100
$ws.Range("D25").Value = 79
$ws.Range("E25").Value = 26.582278481012
$ws.Range("F25").Value = 351
$ws.Range("G25").Value = 337
$ws.Range("H25").Value = 4.154302670623
$ws.Range("I25").Value = 2571
$ws.Range("J25").Value = 2584
$ws.Range("K25").Value = -0.503095975232
$ws.Range("L25").Value = 15.758667266996
$ws.Range("M25").Value = -15.980392156862
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 11
$ws.Range("E26").Value = -72.727272727272
$ws.Range("F26").Value = 13
$ws.Range("G26").Value = 28
$ws.Range("H26").Value = -53.571428571428
$ws.Range("I26").Value = 140
$ws.Range("J26").Value = 173
$ws.Range("K26").Value = -19.07514450867
$ws.Range("L26").Value = 0.719424460431
$ws.Range("C27").Value = 15
$ws.Range("D27").Value = 11
$ws.Range("E27").Value = 36.363636363636
$ws.Range("G27").Value = 54
$ws.Range("H27").Value = 1.851851851851
$ws.Range("I27").Value = 345
$ws.Range("J27").Value = 390
$ws.Range("K27").Value = -11.538461538461
$ws.Range("L27").Value = -0.28901734104
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 7
$ws.Range("E28").Value = -85.714285714285
$ws.Range("F28").Value = 10
$ws.Range("G28").Value = 16
$ws.Range("H28").Value = -37.5
$ws.Range("I28").Value = 84
$ws.Range("J28").Value = 115
$ws.Range("K28").Value = -26.95652173913
$ws.Range("L28").Value = -39.568345323741
$ws.Range("M28").Value = -25.663716814159
$ws.Range("N28").Value = -81.935483870967
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 5
$ws.Range("E29").Value = -80
$ws.Range("F29").Value = 10
$ws.Range("G29").Value = 14
$ws.Range("H29").Value = -28.571428571428
$ws.Range("I29").Value = 76
$ws.Range("J29").Value = 96
$ws.Range("K29").Value = -20.833333333333
$ws.Range("L29").Value = -38.709677419354
$ws.Range("M29").Value = -22.448979591836
$ws.Range("N29").Value = -82.201405152224
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 6
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 50
$ws.Range("I30").Value = 39
$ws.Range("J30").Value = 62
$ws.Range("K30").Value = -37.096774193548
$ws.Range("L30").Value = -13.333333333333

# --- Hate Crimes Week-to-Date 2023 figure becomes a literal "0" (text), matching
#     the zero-count label style used elsewhere in the column (e.g. the N/A dash
#     cells), rather than a plain number. Force it to store as text "0" while
#     reusing the label-style formatting (style of A30) instead of the numeric
#     "#,##0" style that C-column figures normally use. ---
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "0"
$ws.Range("A30").Copy()
$ws.Range("C30").PasteSpecial(-4122)
